$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lpl", "Vldlr", "ECs", 3, 1, 109.3723396666667, 328.117019, 0.3006244632995563, 0.3006244632995563, 3, 1, 0.463191, 1.389573, 0.0353316468093919, 0.0353316468093919, 50.660283382543, 455.942550442887, 0.01062155735956292, 0.01062155735956292),
    @("ECs", "Lpl", "Vldlr", "FAPs", 3, 1, 109.3723396666667, 328.117019, 0.3006244632995563, 0.3006244632995563, 3, 1, 9.932929333333332, 29.798788, 0.7576717833204485, 0.7576717833204486, 1086.38772093033, 9777.489488372972, 0.2277746732179276, 0.2277746732179276),
    @("ECs", "Lpl", "Vldlr", "M2", 3, 1, 109.3723396666667, 328.117019, 0.3006244632995563, 0.3006244632995563, 2, 0.6666666666666666, 0.1285543333333333, 0.385663, 0.009805968382697785, 0.009805968382697785, 14.06028821095522, 126.542593898597, 0.00294791398218094, 0.00294791398218094),
    @("ECs", "Lpl", "Vldlr", "sCs", 3, 1, 109.3723396666667, 328.117019, 0.3006244632995563, 0.3006244632995563, 3, 1, 2.585130333333333, 7.755391, 0.1971906014874617, 0.1971906014874618, 282.7417528999366, 2544.675776099429, 0.05928031873988487, 0.05928031873988489),
    @("FAPs", "Lpl", "Vldlr", "ECs", 3, 1, 170.751104, 512.2533120000001, 0.4693321835689973, 0.4693321835689973, 3, 1, 0.463191, 1.389573, 0.0353316468093919, 0.0353316468093919, 79.09037461286401, 711.813371515776, 0.0165822789461405, 0.0165822789461405),
    @("FAPs", "Lpl", "Vldlr", "FAPs", 3, 1, 170.751104, 512.2533120000001, 0.4693321835689973, 0.4693321835689973, 3, 1, 9.932929333333332, 29.798788, 0.7576717833204485, 0.7576717833204486, 1696.058649620651, 15264.52784658586, 0.3555997524944023, 0.3555997524944023),
    @("FAPs", "Lpl", "Vldlr", "M2", 3, 1, 170.751104, 512.2533120000001, 0.4693321835689973, 0.4693321835689973, 2, 0.6666666666666666, 0.1285543333333333, 0.385663, 0.009805968382697785, 0.009805968382697785, 21.95079434065067, 197.557149065856, 0.0046022565530601, 0.0046022565530601),
    @("FAPs", "Lpl", "Vldlr", "sCs", 3, 1, 170.751104, 512.2533120000001, 0.4693321835689973, 0.4693321835689973, 3, 1, 2.585130333333333, 7.755391, 0.1971906014874617, 0.1971906014874618, 441.4138584005547, 3972.724725604993, 0.09254789557539438, 0.09254789557539439),
    @("M2", "Lpl", "Vldlr", "ECs", 3, 1, 68.96861966666667, 206.905859, 0.1895694499632422, 0.1895694499632422, 3, 1, 0.463191, 1.389573, 0.0353316468093919, 0.0353316468093919, 31.945643912023, 287.510795208207, 0.006697800851951962, 0.006697800851951962),
    @("M2", "Lpl", "Vldlr", "FAPs", 3, 1, 68.96861966666667, 206.905859, 0.1895694499632422, 0.1895694499632422, 3, 1, 9.932929333333332, 29.798788, 0.7576717833204485, 0.7576717833204486, 685.0604253665435, 6165.543828298893, 0.1436314232167262, 0.1436314232167262),
    @("M2", "Lpl", "Vldlr", "M2", 3, 1, 68.96861966666667, 206.905859, 0.1895694499632422, 0.1895694499632422, 2, 0.6666666666666666, 0.1285543333333333, 0.385663, 0.009805968382697785, 0.009805968382697785, 8.866214922168556, 79.795934299517, 0.001858912032664962, 0.001858912032664962),
    @("M2", "Lpl", "Vldlr", "sCs", 3, 1, 68.96861966666667, 206.905859, 0.1895694499632422, 0.1895694499632422, 3, 1, 2.585130333333333, 7.755391, 0.1971906014874617, 0.1971906014874618, 178.2928707484299, 1604.635836735869, 0.037381313861899, 0.037381313861899),
    @("sCs", "Lpl", "Vldlr", "ECs", 3, 1, 14.72510066666667, 44.175302, 0.0404739031682042, 0.04047390316820419, 3, 1, 0.463191, 1.389573, 0.0353316468093919, 0.0353316468093919, 6.820534102894, 61.384806926046, 0.001430009651736519, 0.001430009651736518),
    @("sCs", "Lpl", "Vldlr", "FAPs", 3, 1, 14.72510066666667, 44.175302, 0.0404739031682042, 0.04047390316820419, 3, 1, 9.932929333333332, 29.798788, 0.7576717833204485, 0.7576717833204486, 146.2633843482195, 1316.370459133976, 0.03066593439139242, 0.03066593439139242),
    @("sCs", "Lpl", "Vldlr", "M2", 3, 1, 14.72510066666667, 44.175302, 0.0404739031682042, 0.04047390316820419, 2, 0.6666666666666666, 0.1285543333333333, 0.385663, 0.009805968382697785, 0.009805968382697785, 1.892975499469556, 17.036779495226, 0.0003968858147917821, 0.000396885814791782),
    @("sCs", "Lpl", "Vldlr", "sCs", 3, 1, 14.72510066666667, 44.175302, 0.0404739031682042, 0.04047390316820419, 3, 1, 2.585130333333333, 7.755391, 0.1971906014874617, 0.1971906014874618, 38.06630439478689, 342.596739553082, 0.007981073310283469, 0.007981073310283469)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}
